# BFA - all master data
# Fix the individual_type reference sheet:
#  - is_active column becomes a real boolean (TRUE) instead of the text "TRUE"
#  - the French / Arabic "name" labels are corrected
#  - the ad-hoc bold/border formatting that had crept onto the data rows is cleared
#  - selection is moved to the data range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray per-row formatting (border/alignment/font) that had been
# applied to the data rows A2:D7 so they fall back to the sheet's default style.
$ws.Range("A2:D7").Style = "Normal"

# Re-enter the name values (fra/ara rows) with their corrected text.
$ws.Range("C4").Value = "Ã‰tranger"
$ws.Range("C5").Value = "Non-Ã©tranger"
$ws.Range("C6").Value = "Ø£Ø¬Ù†Ø¨ÙŠ"
$ws.Range("C7").Value = "ØºÙŠØ± Ø£Ø¬Ù†Ø¨ÙŠ"

# is_active: store a native boolean TRUE instead of the text "TRUE".
$ws.Range("D2:D7").Value = $true

# Update the current selection to the data range, as left by the editor.
[void]$ws.Range("A2:D7").Select()
